# Data update: refresh enrollment counts (Pagos/Inscricoes homologadas/Inscritos)
# on the "Inscricoes" sheet for the 2025/1 Superior registration summary table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 34
$ws.Range("H2").Value = 46
$ws.Range("F4").Value = 8
$ws.Range("H4").Value = 22
$ws.Range("E6").Value = 8
$ws.Range("E7").Value = 9
$ws.Range("F8").Value = 10
$ws.Range("H8").Value = 11
$ws.Range("F9").Value = 13
$ws.Range("H9").Value = 22
$ws.Range("E15").Value = 163
$ws.Range("F15").Value = 91
$ws.Range("H15").Value = 132
$ws.Range("E17").Value = 126
$ws.Range("F17").Value = 63
$ws.Range("H17").Value = 95
$ws.Range("F18").Value = 52
$ws.Range("H18").Value = 88
$ws.Range("E19").Value = 61
$ws.Range("F19").Value = 33
$ws.Range("H19").Value = 46
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = 3
$ws.Range("E28").Value = 21
$ws.Range("F28").Value = 18
$ws.Range("H28").Value = 20
$ws.Range("E33").Value = 45
$ws.Range("F35").Value = 4
$ws.Range("H35").Value = 5
$ws.Range("E36").Value = 109
$ws.Range("F36").Value = 49
$ws.Range("H36").Value = 81
$ws.Range("E37").Value = 57
$ws.Range("F37").Value = 33
$ws.Range("H37").Value = 45
$ws.Range("E40").Value = 24
$ws.Range("F40").Value = 13
$ws.Range("H40").Value = 16
$ws.Range("F41").Value = 18
$ws.Range("H41").Value = 29
$ws.Range("F42").Value = 19
$ws.Range("H42").Value = 28
$ws.Range("F44").Value = 15
$ws.Range("H44").Value = 25
$ws.Range("F46").Value = 11
$ws.Range("H46").Value = 20
$ws.Range("F47").Value = 39
$ws.Range("H47").Value = 50
$ws.Range("F51").Value = 11
$ws.Range("H51").Value = 11
$ws.Range("F57").Value = 5
$ws.Range("H57").Value = 9
$ws.Range("E60").Value = 19
$ws.Range("E62").Value = 49
$ws.Range("F62").Value = 13
$ws.Range("H62").Value = 27
$ws.Range("E66").Value = 35
$ws.Range("F66").Value = 23
$ws.Range("H66").Value = 31
$ws.Range("F68").Value = 11
$ws.Range("H68").Value = 15
$ws.Range("E70").Value = 45
$ws.Range("F72").Value = 25
$ws.Range("H72").Value = 36
$ws.Range("F74").Value = 8
$ws.Range("H74").Value = 12
$ws.Range("E79").Value = 40
$ws.Range("F79").Value = 19
$ws.Range("H79").Value = 31
$ws.Range("F80").Value = 13
$ws.Range("H80").Value = 25
$ws.Range("F87").Value = 5
$ws.Range("H87").Value = 12
$ws.Range("E88").Value = 26
$ws.Range("F88").Value = 14
$ws.Range("H88").Value = 22
$ws.Range("E89").Value = 47
$ws.Range("F89").Value = 21
$ws.Range("H89").Value = 28
